# Add filter criteria for household table when filtering household_members
#
# 1. Insert a new column before column I on the "survey" sheet, inheriting
#    the width of the column to its left (G:H).
# 2. Populate the new column's header (I1) and the new filter-criteria cells
#    in row 7 (G7, H7, I7).
# 3. Make "survey" the active sheet/tab (was "settings"), with the viewport
#    scrolled so column E is at the left edge and G7 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- Insert a new column before I, matching the width of column H ---
$leftWidth = $ws.Columns("H:H").ColumnWidth
$ws.Columns("I:I").Insert()
$ws.Columns("I:I").ColumnWidth = $leftWidth

# --- New values, written in the order that matches the target shared
#     string table (household_id = ?, selectionArgs.cell_type, formula,
#     [ data('household_id') ]) ---
$ws.Range("G7").Value = "household_id = ?"
$ws.Range("I1").Value = "selectionArgs.cell_type"
$ws.Range("I7").Value = "formula"
$ws.Range("H7").Value = "[ data('household_id') ]"

# --- Make survey the active sheet, scrolled/selected as in the edit ---
$ws.Activate()
$ws.Range("G7").Select()
